$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G36").Value = 0.8269467468186171
$ws.Range("H36").Value = 15.12096357345581
$ws.Range("J36").Value = 15.12096357345581
$ws.Range("K36").Value = 0

$ws.Range("G37").Value = 0.8346286309244421
$ws.Range("H37").Value = 0.854478120803833
$ws.Range("J37").Value = 0.854478120803833
$ws.Range("K37").Value = 1

$ws.Range("G38").Value = 0.8258022611941038
$ws.Range("H38").Value = 59.09400129318237
$ws.Range("J38").Value = 0.5909400129318237
$ws.Range("K38").Value = 2

$ws.Range("G39").Value = 0.8278452114356494
$ws.Range("H39").Value = 5.610887050628662
$ws.Range("J39").Value = 5.610887050628662
$ws.Range("K39").Value = 3

$ws.Range("G40").Value = 0.8294518917544784
$ws.Range("H40").Value = 7.493222951889038
$ws.Range("J40").Value = 7.493222951889038
$ws.Range("K40").Value = 4

$ws.Range("G41").Value = 0.8408885011309946
$ws.Range("H41").Value = 1.294256448745728
$ws.Range("J41").Value = 1.294256448745728
$ws.Range("K41").Value = 5

$ws.Range("G42").Value = 0.8410602263105685
$ws.Range("H42").Value = 1.466903209686279
$ws.Range("J42").Value = 1.466903209686279
$ws.Range("K42").Value = 6

$ws.Range("G43").Value = 0.8412392467969956
$ws.Range("H43").Value = 27.02620673179626
$ws.Range("J43").Value = 0.2702620673179627
$ws.Range("K43").Value = 7

$ws.Range("G44").Value = 0.8425120240930364
$ws.Range("H44").Value = 29.18382358551025
$ws.Range("J44").Value = 0.2918382358551025
$ws.Range("K44").Value = 8

$ws.Range("G45").Value = 0.8265007921117307
$ws.Range("H45").Value = 5.550853252410889
$ws.Range("J45").Value = 5.550853252410889
$ws.Range("K45").Value = 9

$ws.Range("G46").Value = 0.8275464739064742
$ws.Range("H46").Value = 7.196524858474731
$ws.Range("J46").Value = 7.196524858474731
$ws.Range("K46").Value = 10

$ws.Range("G47").Value = 0.8265007921117307
$ws.Range("H47").Value = 0.2990939617156982
$ws.Range("J47").Value = 0.2990939617156982
$ws.Range("K47").Value = 11

$ws.Range("G48").Value = 0.8354558235653933
$ws.Range("H48").Value = 0.8191165924072266
$ws.Range("J48").Value = 0.8191165924072266
$ws.Range("K48").Value = 12

$ws.Range("G49").Value = 0.8400244619039816
$ws.Range("H49").Value = 24.14540815353394
$ws.Range("J49").Value = 0.2414540815353393
$ws.Range("K49").Value = 13

$ws.Range("G50").Value = 0.8417985642286322
$ws.Range("H50").Value = 28.50107908248901
$ws.Range("J50").Value = 0.2850107908248901
$ws.Range("K50").Value = 14

$ws.Range("G51").Value = 0.8240171630446314
$ws.Range("H51").Value = 6.64790940284729
$ws.Range("J51").Value = 6.64790940284729
$ws.Range("K51").Value = 15

$ws.Range("G52").Value = 0.8233166127014845
$ws.Range("H52").Value = 7.618446826934814
$ws.Range("J52").Value = 7.618446826934814
$ws.Range("K52").Value = 16

$ws.Range("G53").Value = 0.8259885769873838
$ws.Range("H53").Value = 0.3562424182891846
$ws.Range("J53").Value = 0.3562424182891846
$ws.Range("K53").Value = 17

$ws.Range("G54").Value = 0.8352137673985931
$ws.Range("H54").Value = 0.6476178169250488
$ws.Range("J54").Value = 0.6476178169250488
$ws.Range("K54").Value = 18

$ws.Range("G55").Value = 0.8419818878888116
$ws.Range("H55").Value = 24.97894978523254
$ws.Range("J55").Value = 0.2497894978523254
$ws.Range("K55").Value = 19

$ws.Range("G56").Value = 0.8395874835963968
$ws.Range("H56").Value = 26.64813280105591
$ws.Range("J56").Value = 0.2664813280105591
$ws.Range("K56").Value = 20

$ws.Range("G57").Value = 0.8270412157559089
$ws.Range("H57").Value = 5.440364837646484
$ws.Range("J57").Value = 5.440364837646484
$ws.Range("K57").Value = 21

$ws.Range("G58").Value = 0.8301606783933332
$ws.Range("H58").Value = 7.428974866867065
$ws.Range("J58").Value = 7.428974866867065
$ws.Range("K58").Value = 22

$ws.Range("G59").Value = 0.8392204662242412
$ws.Range("H59").Value = 0.5779387950897217
$ws.Range("J59").Value = 0.5779387950897217
$ws.Range("K59").Value = 23

$ws.Range("G60").Value = 0.8343136216891123
$ws.Range("H60").Value = 0.5903444290161133
$ws.Range("J60").Value = 0.5903444290161133
$ws.Range("K60").Value = 24

$ws.Range("G61").Value = 0.8390526741666228
$ws.Range("H61").Value = 20.30105566978455
$ws.Range("J61").Value = 0.2030105566978455
$ws.Range("K61").Value = 25

$ws.Range("G62").Value = 0.838786480049757
$ws.Range("H62").Value = 25.31743574142456
$ws.Range("J62").Value = 0.2531743574142456
$ws.Range("K62").Value = 26

$ws.Range("G63").Value = 0.8427752260768016
$ws.Range("H63").Value = 55.05057048797607
$ws.Range("J63").Value = 55.05057048797607
$ws.Range("K63").Value = 27

$ws.Range("G64").Value = 0.8427752260768016
$ws.Range("H64").Value = 60.56906414031982
$ws.Range("J64").Value = 60.56906414031982
$ws.Range("K64").Value = 28

$ws.Range("G65").Value = 0.842214597604513
$ws.Range("H65").Value = 330.8678059577942
$ws.Range("J65").Value = 33.08678059577942
$ws.Range("K65").Value = 29

$ws.Range("G66").Value = 0.842214597604513
$ws.Range("H66").Value = 407.5261158943176
$ws.Range("J66").Value = 40.75261158943177
$ws.Range("K66").Value = 30
